$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 4167672.8
$ws.Range("I33").Value = 5555773
$ws.Range("K33").Value = 5555773
$ws.Range("M33").Value = -5555544

# Row 82
$ws.Range("H82").Value = 1323.6666
$ws.Range("I82").Value = 1323.6666
$ws.Range("K82").Value = 3970.9998
$ws.Range("M82").Value = -3564.9998

# Row 85
$ws.Range("H85").Value = 1323.6666
$ws.Range("I85").Value = 1323.6666
$ws.Range("K85").Value = 3970.9998
$ws.Range("M85").Value = -2566.9998

# Row 113
$ws.Range("H113").Value = 72976.875
$ws.Range("J113").Value = 12994.8
$ws.Range("L113").Value = 12994.8
$ws.Range("N113").Value = -19502.8

# Row 127
$ws.Range("J127").Value = 2000
$ws.Range("L127").Value = 6000
$ws.Range("N127").Value = -15920

# Row 129
$ws.Range("H129").Value = 1686.1
$ws.Range("I129").Value = 1686.1
$ws.Range("K129").Value = 5058.299999999999
$ws.Range("M129").Value = -58.29999999999927

# Row 131
$ws.Range("H131").Value = 2550
$ws.Range("I131").Value = 2550
$ws.Range("K131").Value = 7650
$ws.Range("M131").Value = -2610

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3447.4478
$ws.Range("I32").Value = 3019.0322
$ws.Range("K32").Value = 3019.0322
$ws.Range("M32").Value = -2732.0322

# Row 45
$ws.Range("H45").Value = 5474.147
$ws.Range("I45").Value = 8708.6
$ws.Range("K45").Value = 8708.6
$ws.Range("M45").Value = -8331.6

# Row 74
$ws.Range("H74").Value = 7806.8887
$ws.Range("I74").Value = 925.44446
$ws.Range("K74").Value = 925.44446
$ws.Range("M74").Value = -51.44446000000005

# Row 77
$ws.Range("H77").Value = 7806.8887
$ws.Range("I77").Value = 925.44446
$ws.Range("K77").Value = 4627.2223
$ws.Range("M77").Value = -259.2223000000004

# Row 110
$ws.Range("H110").Value = 5730.5674
$ws.Range("I110").Value = 6147.7617
$ws.Range("K110").Value = 6147.7617
$ws.Range("M110").Value = -4102.7617

# Row 132
$ws.Range("H132").Value = 3331.375
$ws.Range("I132").Value = 2693.5
$ws.Range("K132").Value = 8080.5
$ws.Range("M132").Value = -5550.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2544.111
$ws.Range("I86").Value = 2452.4666
$ws.Range("K86").Value = 2452.4666
$ws.Range("M86").Value = -1329.4666

# Row 89
$ws.Range("H89").Value = 2544.111
$ws.Range("I89").Value = 2452.4666
$ws.Range("K89").Value = 12262.333
$ws.Range("M89").Value = -6646.333000000001

# Row 139
$ws.Range("H139").Value = 75296.2
$ws.Range("J139").Value = 74120.25
$ws.Range("L139").Value = 74120.25
$ws.Range("N139").Value = -84400.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1503
$ws.Range("I16").Value = 1586.8334
$ws.Range("K16").Value = 1586.8334
$ws.Range("M16").Value = -1299.8334

# Row 31
$ws.Range("H31").Value = 38115.035
$ws.Range("I31").Value = 49014.332
$ws.Range("J31").Value = 9504.375
$ws.Range("K31").Value = 49014.332
$ws.Range("L31").Value = 9504.375
$ws.Range("M31").Value = -48719.332
$ws.Range("N31").Value = -10094.375

# Row 34
$ws.Range("H34").Value = 38115.035
$ws.Range("I34").Value = 49014.332
$ws.Range("J34").Value = 9504.375
$ws.Range("K34").Value = 49014.332
$ws.Range("L34").Value = 9504.375
$ws.Range("M34").Value = -48812.332
$ws.Range("N34").Value = -9908.375

# Row 58
$ws.Range("H58").Value = 2881.6667
$ws.Range("I58").Value = 2941.923
$ws.Range("K58").Value = 2941.923
$ws.Range("M58").Value = -2738.923

# Row 62
$ws.Range("H62").Value = 15000
$ws.Range("I62").Value = 15000
$ws.Range("K62").Value = 15000
$ws.Range("M62").Value = -14376

# Row 65
$ws.Range("H65").Value = 15000
$ws.Range("I65").Value = 15000
$ws.Range("K65").Value = 75000
$ws.Range("M65").Value = -71880

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 107
$ws.Range("H107").Value = 308.2381
$ws.Range("I107").Value = 323.9375
$ws.Range("J107").Value = 258
$ws.Range("K107").Value = 323.9375
$ws.Range("L107").Value = 258
$ws.Range("M107").Value = 1596.0625
$ws.Range("N107").Value = -4098

# Row 113
$ws.Range("H113").Value = 1503
$ws.Range("I113").Value = 1586.8334
$ws.Range("K113").Value = 1586.8334
$ws.Range("M113").Value = 583.1666

# Row 136
$ws.Range("H136").Value = 2881.6667
$ws.Range("I136").Value = 2941.923
$ws.Range("K136").Value = 8825.769
$ws.Range("M136").Value = -6275.769

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 307.89285
$ws.Range("I23").Value = 191
$ws.Range("J23").Value = 600.125
$ws.Range("K23").Value = 573
$ws.Range("L23").Value = 1800.375
$ws.Range("M23").Value = -338
$ws.Range("N23").Value = -2270.375

# Row 122
$ws.Range("H122").Value = 3014.4285
$ws.Range("I122").Value = 299.5
$ws.Range("J122").Value = 4100.4
$ws.Range("K122").Value = 2695.5
$ws.Range("L122").Value = 36903.6
$ws.Range("M122").Value = -245.5
$ws.Range("N122").Value = -41803.6

# Row 132
$ws.Range("H132").Value = 1459.8667
$ws.Range("I132").Value = 989.8
$ws.Range("K132").Value = 8908.199999999999
$ws.Range("M132").Value = -6378.199999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 15000275
$ws.Range("I11").Value = 17142886
$ws.Range("K11").Value = 17142886
$ws.Range("M11").Value = -17142747

# Row 14
$ws.Range("H14").Value = 2600600.8
$ws.Range("I14").Value = 3714786.2
$ws.Range("J14").Value = 834.6667
$ws.Range("K14").Value = 3714786.2
$ws.Range("L14").Value = 834.6667
$ws.Range("M14").Value = -3714618.2
$ws.Range("N14").Value = -1170.6667

# Row 70
$ws.Range("H70").Value = 16190
$ws.Range("I70").Value = 10554.5
$ws.Range("K70").Value = 10554.5
$ws.Range("M70").Value = -10284.5

# Row 73
$ws.Range("H73").Value = 16190
$ws.Range("I73").Value = 10554.5
$ws.Range("K73").Value = 10554.5
$ws.Range("M73").Value = -9618.5

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 45500
$ws.Range("J3").Value = 45500
$ws.Range("L3").Value = 45500
$ws.Range("N3").Value = -45724

# Row 15
$ws.Range("H15").Value = 45500
$ws.Range("J15").Value = 45500
$ws.Range("L15").Value = 45500
$ws.Range("N15").Value = -45840

# Row 68
$ws.Range("H68").Value = 5763.125
$ws.Range("I68").Value = 4015
$ws.Range("K68").Value = 4015
$ws.Range("M68").Value = -3266

# Row 71
$ws.Range("H71").Value = 5763.125
$ws.Range("I71").Value = 4015
$ws.Range("K71").Value = 20075
$ws.Range("M71").Value = -16331

# Row 93
$ws.Range("H93").Value = 2444.3684
$ws.Range("I93").Value = 2638
$ws.Range("J93").Value = 798.5
$ws.Range("K93").Value = 2638
$ws.Range("L93").Value = 798.5
$ws.Range("M93").Value = -1390
$ws.Range("N93").Value = -3294.5

# Row 132
$ws.Range("H132").Value = 3511.1353
$ws.Range("I132").Value = 2803.5862
$ws.Range("J132").Value = 6076
$ws.Range("K132").Value = 8410.758600000001
$ws.Range("L132").Value = 18228
$ws.Range("M132").Value = -5880.758600000001
$ws.Range("N132").Value = -23288

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 505000
$ws.Range("I3").Value = 505000
$ws.Range("K3").Value = 505000
$ws.Range("M3").Value = -504886

# Row 14
$ws.Range("H14").Value = 740666.3
$ws.Range("I14").Value = 793499.7
$ws.Range("J14").Value = 999
$ws.Range("K14").Value = 793499.7
$ws.Range("L14").Value = 999
$ws.Range("M14").Value = -793331.7
$ws.Range("N14").Value = -1335

# Row 74
$ws.Range("H74").Value = 11690.429
$ws.Range("J74").Value = 11690.429
$ws.Range("L74").Value = 11690.429
$ws.Range("N74").Value = -13562.429

# Row 77
$ws.Range("H77").Value = 11690.429
$ws.Range("J77").Value = 11690.429
$ws.Range("L77").Value = 35071.287
$ws.Range("N77").Value = -44431.287

# Row 132
$ws.Range("H132").Value = 3327.0938
$ws.Range("I132").Value = 3272.1072
$ws.Range("K132").Value = 9816.321599999999
$ws.Range("M132").Value = -7286.321599999999
